# ajuste no formato das coordenadas e mapeado os pontos do 14 ao 32
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: BFS
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BFS")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = "[4, 5]"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1.333333333333333
$ws.Cells.Item(2, 7).Value = 0.000301361083984375

$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "[3, 6, 9, 13, 17, 20, 24]"
$ws.Cells.Item(3, 5).Value = 27
$ws.Cells.Item(3, 6).Value = 1.037037037037037
$ws.Cells.Item(3, 7).Value = 0.000102996826171875

$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = "[11, 12, 13, 17, 20]"
$ws.Cells.Item(4, 5).Value = 22
$ws.Cells.Item(4, 6).Value = 1.136363636363636
$ws.Cells.Item(4, 7).Value = 0.00008034706115722656

# ---------------------------------------------------------------
# Sheet: DFS
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DFS")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = "[4, 1, 2, 3, 6, 5]"
$ws.Cells.Item(2, 7).Value = 0.0001111030578613281

$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "[3, 2, 5, 4, 7, 8, 9, 6, 27, 28, 29, 13, 12, 11, 10, 14, 15, 16, 17, 20, 19, 18, 21, 22, 23, 24]"
$ws.Cells.Item(3, 5).Value = 29
$ws.Cells.Item(3, 6).Value = 0.03448275862068965
$ws.Cells.Item(3, 7).Value = 0.00007557868957519531

$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = "[11, 10, 7, 4, 1, 2, 3, 6, 5, 8, 9, 28, 29, 13, 12, 16, 15, 14, 18, 19, 20]"
$ws.Cells.Item(4, 5).Value = 24
$ws.Cells.Item(4, 6).Value = 0.08333333333333333
$ws.Cells.Item(4, 7).Value = 0.00007605552673339844

# ---------------------------------------------------------------
# Sheet: BCU
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BCU")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = "[4, 5]"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 0.4
$ws.Cells.Item(2, 7).Value = 0.0000820159912109375

$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "[3, 25, 26, 27, 28, 29, 13, 17, 20, 24]"
$ws.Cells.Item(3, 5).Value = 23
$ws.Cells.Item(3, 6).Value = 0.6285714285714286
$ws.Cells.Item(3, 7).Value = 0.00009298324584960938

$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = "[11, 15, 16, 17, 20]"
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(4, 6).Value = 0.5263157894736842
$ws.Cells.Item(4, 7).Value = 0.00007653236389160156

# ---------------------------------------------------------------
# Sheet: A_Estrela_Euclidiano
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("A_Estrela_Euclidiano")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = "[4, 5]"
$ws.Cells.Item(2, 6).Value = 3.0625
$ws.Cells.Item(2, 7).Value = 0.00009417533874511719

$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "[3, 25, 26, 27, 28, 29, 30, 31, 20, 24]"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 3.0625
$ws.Cells.Item(3, 7).Value = 0.000095367431640625

$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = "[11, 15, 16, 17, 20]"
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 3.0625
$ws.Cells.Item(4, 7).Value = 0.000110626220703125

# ---------------------------------------------------------------
# Sheet: A_Estrela_Haversiano
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("A_Estrela_Haversiano")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = "[4, 5]"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 3.0625
$ws.Cells.Item(2, 7).Value = 0.0002543926239013672

$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "[3, 25, 26, 27, 28, 29, 30, 31, 32, 24]"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 3.0625
$ws.Cells.Item(3, 7).Value = 0.0001401901245117188

$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = "[11, 15, 16, 19, 20]"
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 3.0625
$ws.Cells.Item(4, 7).Value = 0.0001769065856933594
